$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update activation date (row 8: B8/C8) - force text so Excel doesn't
# auto-convert the date-like string into a date serial value.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2022"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "01/01/2022"

# Update "Programa:" detailed content (row 16: B16/C16) - Portuguese syllabus bullet list
$programaText = "- Reologia de fluidos,- Dimensionamento de tubulações,- Acessórios e bombeamento para fluidos industriais,- Agitação e mistura,- Caracterização de partículas e leito de partículas,- Sedimentação,- Filtração,- Processos com membranas.- Operações unitárias de troca térmica: trocadores de calor e evaporadores."
$ws.Range("B16").Value = $programaText
$ws.Range("C16").Value = $programaText

# Update "Syllabus:" detailed content (row 17: B17/C17) - English syllabus bullet list
$syllabusText = "- fluid rheology,- Sizing of pipes,- Accessories and pumping for industrial fluids,- Stirring and mixing,- Particle characterization and particle bed,- Sedimentation,- Filtration,- Processes with membranes.- Unit heat exchange operations: heat exchangers and evaporators."
$ws.Range("B17").Value = $syllabusText
$ws.Range("C17").Value = $syllabusText
